# semana 31 de 2025
# Adds week-31 data: a new column AH ("31") with counts per provider row,
# and a new provider row (SALUD PYP SAS) inserted before the existing
# row for cod_pre 6600102402.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new data row (becomes row 39), pushing the old
#        rows 39-57 down to 40-58 ---------------------------------------
$ws.Rows(39).Insert()

$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "6600102288"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "01"
$ws.Range("C39").Value = "SALUD PYP SAS"

# --- 2. Add the new "31" column header in AH1, matching the bold /
#        centered style used by the rest of the header row ----------------
$ws.Range("AH1").Font.Bold = $true
$ws.Range("AH1").HorizontalAlignment = -4108
$ws.Range("AH1").NumberFormat = "@"
$ws.Range("AH1").Value = "31"

# --- 3. Fill in the week-31 (AH) values for every data row ----------------
$ws.Range("AH2").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AH6").Value = 26
$ws.Range("AH7").Value = 3
$ws.Range("AH8").Value = 29
$ws.Range("AH10").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AH12").Value = 0
$ws.Range("AH13").Value = 0
$ws.Range("AH14").Value = 0
$ws.Range("AH15").Value = 0
$ws.Range("AH16").Value = 0
$ws.Range("AH17").Value = 0
$ws.Range("AH23").Value = 0
$ws.Range("AH24").Value = 0
$ws.Range("AH25").Value = 2
$ws.Range("AH28").Value = 55
$ws.Range("AH29").Value = 4
$ws.Range("AH30").Value = 15
$ws.Range("AH31").Value = 0
$ws.Range("AH32").Value = 0
$ws.Range("AH34").Value = 0
$ws.Range("AH35").Value = 17
$ws.Range("AH36").Value = 1
$ws.Range("AH37").Value = 0
$ws.Range("AH38").Value = 0
$ws.Range("AH39").Value = 0
$ws.Range("AH40").Value = 0
$ws.Range("AH41").Value = 0
$ws.Range("AH42").Value = 0
$ws.Range("AH43").Value = 0
$ws.Range("AH45").Value = 0
$ws.Range("AH46").Value = 0
$ws.Range("AH47").Value = 0
$ws.Range("AH48").Value = 0
$ws.Range("AH49").Value = 0
$ws.Range("AH50").Value = 0
$ws.Range("AH51").Value = 0
$ws.Range("AH53").Value = 0
$ws.Range("AH54").Value = 0
$ws.Range("AH55").Value = 0
$ws.Range("AH56").Value = 0
$ws.Range("AH57").Value = 0
$ws.Range("AH58").Value = 0
